# TC42_Canine_Filter_Breed-ShihTzu.xlsx - "updated ICDC study obj ids 09/02"
#
# The Cypher query stored in cell B2 of the "startup" sheet is extended so
# the WITH DISTINCT clause also threads `demo.weight` through (it is already
# consumed further down in the RETURN clause's weight CASE expression).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$cell = $ws.Range("B2")
$query = $cell.Value2

$oldClause = "WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age`n"
$newClause = "WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight`n"

$query = $query.Replace($oldClause, $newClause)
$cell.Value2 = $query

# Reflect the editing session's final selection on the sheet.
$ws.Range("C4").Select() | Out-Null
